$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price / volume(1h) snapshot pulled by the scheduled scraper.
# Plain decimal "Price" values must stay TEXT (the sheet stores e.g. "66.654.70"
# as a literal string), so force Text format before writing those so Excel
# does not auto-convert them to numbers; percentage cells already contain
# non-numeric characters (%, padding spaces) and are safe to set directly.

$ws.Range('D2').Value = '66.654.70'
$ws.Range('E2').Value = '  +0.23%  '
$ws.Range('D3').Value = '3.252.56'
$ws.Range('E3').Value = '  +2.14%  '
$ws.Range('E4').Value = '  +0.17%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '604.13'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.16%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '157.91'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +1.28%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = '3.253.67'
$ws.Range('E8').Value = '  +2.28%  '
$ws.Range('E9').Value = '  -0.17%  '
$ws.Range('E10').Value = '  +1.93%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '5.92'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +3.99%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.507'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -1.17%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.0000271'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +1.28%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '39.39'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +1.25%  '
$ws.Range('D15').Value = '3.792.48'
$ws.Range('E15').Value = '  +2.31%  '
$ws.Range('D16').Value = '66.740.08'
$ws.Range('E16').Value = '  +0.37%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '7.38'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -0.76%  '
$ws.Range('D18').Value = '3.258.50'
$ws.Range('E18').Value = '  +2.62%  '
$ws.Range('E19').Value = '  +1.27%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '508.34'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -1.11%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '15.39'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -1.00%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '0.753'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +2.35%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '8.06'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -2.58%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '14.82'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -1.08%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '86.49'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +2.08%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '0.164'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +82.97%  '
$ws.Range('E27').Value = '  +0.11%  '
$ws.Range('E28').Value = '  +0.18%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '9.09'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -1.62%  '
$ws.Range('E30').Value = '  -1.20%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '6.89'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -2.94%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '2.86'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -8.49%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '28.18'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +0.05%  '
$ws.Range('E34').Value = '  +0.34%  '
$ws.Range('E35').Value = '  -4.60%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '6.43'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -2.04%  '
$ws.Range('D37').Value = '0.0₃0806'
$ws.Range('E37').Value = '  +17.43%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '55.48'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +1.35%  '
$ws.Range('B39').Value = 'dogwifhat'
$ws.Range('C39').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '3.33'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +16.87%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '495.27'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -3.55%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.0428'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +1.07%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.128'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +1.57%  '
$ws.Range('E43').Value = '  -2.72%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.296'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -2.02%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '2.48'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +1.34%  '
$ws.Range('D46').Value = '2.942.74'
$ws.Range('E46').Value = '  +3.16%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '28.58'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +0.24%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '2.45'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +1.39%  '
$ws.Range('E49').Value = '  +1.85%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '2.54'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -2.52%  '
